# "test P7 with -10 percent" -- replace the solver-result values across all
# sheets of the workbook with the re-run (P7, -10%) figures, and drop the
# now-unused 3rd/4th scenario rows (s/index values 2 and 4 for j=5) from the
# y / rho / alpha sheets.

$wb = $excel.ActiveWorkbook

# ---- general: summary KPIs -------------------------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value  = 433.410743643259
$ws.Range("B4").Value  = 0.01999998092651367
$ws.Range("B6").Value  = 33.93074364325894
$ws.Range("B7").Value  = 2.580484450641003
$ws.Range("B8").Value  = 2.580484450641003
$ws.Range("B9").Value  = 235.1
$ws.Range("B10").Value = 164.38

# ---- x: assignment variable ------------------------------------------------
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value  = 1
$ws.Range("B3").Value  = 5
$ws.Range("B4").Value  = 9
$ws.Range("B5").Value  = 7
$ws.Range("B8").Value  = 13
$ws.Range("B9").Value  = 6
$ws.Range("B11").Value = 10
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 8
$ws.Range("B14").Value = 2

# ---- U ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Range("B2").Value  = 1
$ws.Range("B11").Value = 3

# ---- TBar ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value  = 5.468996656383652
$ws.Range("B4").Value  = 36.95051889297534
$ws.Range("B5").Value  = 34.16886835983306
$ws.Range("B6").Value  = 34.76592070603971
$ws.Range("B7").Value  = 30
$ws.Range("B8").Value  = 33.8200130889828
$ws.Range("B10").Value = 32.45367071955468
$ws.Range("B11").Value = 30
$ws.Range("B12").Value = 32.01159140980466
$ws.Range("B13").Value = 36.48759645946009
$ws.Range("B14").Value = 39.53100334361635
$ws.Range("B15").Value = 38.28184163802894

# ---- y: drop the two rows for j=5,s=13 (C=2 and C=4) -----------------------
$ws = $wb.Worksheets.Item("y")
$ws.Rows("7:8").Delete()

# ---- Q: full re-run of the Q table, rows 7..71 ------------------------------
$ws = $wb.Worksheets.Item("Q")
$QVals = @(97.1700000000008, 99.2700000000008, 100.1150000000008, 99.0400000000008, 97.9800000000008, 319.67, 323.35, 324.535, 329.9, 320.095, 224.1799999999995, 224.6649999999995, 201.1149999999995, 218.9699999999995, 207.1049999999995, 226.0399999999994, 247.1799999999994, 221.8549999999994, 238.4549999999994, 224.4749999999995, 69.88000000000022, 76.32000000000022, 72.41000000000021, 72.64500000000022, 72.89000000000021, 177.26, 188.5600000000001, 166.2, 181.48, 167.6700000000001, 180.2450000000007, 183.9900000000007, 178.0900000000007, 188.8100000000007, 179.8350000000007, 153.4099999999999, 167.1249999999999, 139.5349999999999, 154.5, 143.6599999999999, 85.48500000000051, 87.9650000000005, 79.71500000000052, 90.7300000000005, 84.73000000000052, 238.195, 242.67, 239.82, 249.465, 232.75, 224.1799999999995, 224.6649999999995, 201.1149999999995, 218.9699999999995, 207.1049999999995, 319.67, 323.35, 324.535, 329.9, 320.095, 226.0399999999994, 247.1799999999994, 221.8549999999994, 238.4549999999994, 224.4749999999995)
$row = 7
foreach ($v in $QVals) {
    $ws.Cells.Item($row, 3).Value = $v
    $row = $row + 1
}

# ---- R: subset of rows change (j=12 and j=13 blocks) ------------------------
$ws = $wb.Worksheets.Item("R")
$RRows = @(7, 8, 9, 10, 11, 13, 15)
$RVals = @(19.67, 23.35, 24.535, 29.9, 20.095, 0, 0)
for ($i = 0; $i -lt $RRows.Length; $i++) {
    $ws.Cells.Item($RRows[$i], 3).Value = $RVals[$i]
}

# ---- L: j=1 block gains values, j=10 block drops to zero ---------------------
$ws = $wb.Worksheets.Item("L")
$LRows = @(2, 3, 4, 5, 6, 47, 48, 49, 50, 51)
$LVals = @(12.775, 18.32, 13.155, 20.115, 17.825, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $LRows.Length; $i++) {
    $ws.Cells.Item($LRows[$i], 3).Value = $LVals[$i]
}

# ---- rho / alpha: drop the two rows for j=5,s=2 and j=5,s=4 ------------------
$ws = $wb.Worksheets.Item("rho")
$ws.Rows("7:8").Delete()

$ws = $wb.Worksheets.Item("alpha")
$ws.Rows("7:8").Delete()

Write-Host "edits applied"
